# Applies the "additional scraping" edit described in the commit:
#  - Insert a new "Player Info" worksheet as the first sheet, with
#    ID / NAME / BATTING_HAND / BOWL_STYLE columns + one data row for the
#    player (id 4796).
#  - On "ODI Batting": rename column header MATCH_CARD_LINK -> MATCH_CODE
#    and replace the scorecard URLs with just the numeric match code.
#  - On "ODI Bowling": same header + value change for its MATCH_CARD_LINK
#    column (column B there).

$wb = $excel.ActiveWorkbook

$battingSheetRef = $wb.Worksheets.Item("ODI Batting")

# --- 1. Add the new "Player Info" sheet, placed before "ODI Batting" ---
# NOTE: adding a sheet "Before" an existing one repoints the existing
# object reference at the newly-inserted sheet, so the original sheets
# must be re-fetched by name afterwards rather than reused here.
$playerInfo = $wb.Worksheets.Add($battingSheetRef)
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Reuse the exact header-row formatting already used on the other sheets
# (bold, boxed, centered / top-aligned) by copying it across, then fill
# in the header text.
$battingSheet.Range("A1:D1").Copy() | Out-Null
$playerInfo.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# ID is stored as text (not a number) in the source data. Force the cell
# to text format so the numeric-looking value isn't auto-converted to a
# number, then clear the formatting again (leaving the value as text)
# so the cell doesn't end up with a lingering "text format" style that
# the source data never had.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4796"
$playerInfo.Range("A2").ClearFormats()
$playerInfo.Range("B2").Value = "Mitchell Joseph Swepson"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D) ---
# MATCH_CODE values are text (plain match-code numbers extracted from the
# old URL), so force text format before assigning them too, then clear
# the formatting so no style attribute lingers on the cell.
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2:D4").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4564"
$battingSheet.Range("D3").Value = "4565"
$battingSheet.Range("D4").Value = "4597"
$battingSheet.Range("D2:D4").ClearFormats()

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B) ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2:B4").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4564"
$bowlingSheet.Range("B3").Value = "4565"
$bowlingSheet.Range("B4").Value = "4597"
$bowlingSheet.Range("B2:B4").ClearFormats()
